$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col11a1"
$ws.Cells.Item(2, 3).Value = "Ddr1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.05596366666666666
$ws.Cells.Item(2, 8).Value = 0.167891
$ws.Cells.Item(2, 9).Value = 0.1046634785920854
$ws.Cells.Item(2, 10).Value = 0.1046634785920854
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.2113696666666667
$ws.Cells.Item(2, 14).Value = 0.634109
$ws.Cells.Item(2, 15).Value = 0.03795977003925348
$ws.Cells.Item(2, 16).Value = 0.03795977003925347
$ws.Cells.Item(2, 17).Value = 0.01182902156877778
$ws.Cells.Item(2, 18).Value = 0.106461194119
$ws.Cells.Item(2, 19).Value = 0.003973001578863892
$ws.Cells.Item(2, 20).Value = 0.003973001578863891

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col11a1"
$ws.Cells.Item(3, 3).Value = "Ddr1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.05596366666666666
$ws.Cells.Item(3, 8).Value = 0.167891
$ws.Cells.Item(3, 9).Value = 0.1046634785920854
$ws.Cells.Item(3, 10).Value = 0.1046634785920854
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.589504333333333
$ws.Cells.Item(3, 14).Value = 4.768513
$ws.Cells.Item(3, 15).Value = 0.2854582680725092
$ws.Cells.Item(3, 16).Value = 0.2854582680725091
$ws.Cells.Item(3, 17).Value = 0.08895449067588888
$ws.Cells.Item(3, 18).Value = 0.800590416083
$ws.Cells.Item(3, 19).Value = 0.02987705532934084
$ws.Cells.Item(3, 20).Value = 0.02987705532934084

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col11a1"
$ws.Cells.Item(4, 3).Value = "Ddr1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.05596366666666666
$ws.Cells.Item(4, 8).Value = 0.167891
$ws.Cells.Item(4, 9).Value = 0.1046634785920854
$ws.Cells.Item(4, 10).Value = 0.1046634785920854
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.767380666666666
$ws.Cells.Item(4, 14).Value = 11.302142
$ws.Cells.Item(4, 15).Value = 0.6765819618882374
$ws.Cells.Item(4, 16).Value = 0.6765819618882374
$ws.Cells.Item(4, 17).Value = 0.2108364358357777
$ws.Cells.Item(4, 18).Value = 1.897527922522
$ws.Cells.Item(4, 19).Value = 0.07081342168388068
$ws.Cells.Item(4, 20).Value = 0.07081342168388069

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Col11a1"
$ws.Cells.Item(5, 3).Value = "Ddr1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.1337306666666667
$ws.Cells.Item(5, 8).Value = 0.401192
$ws.Cells.Item(5, 9).Value = 0.2501036404769519
$ws.Cells.Item(5, 10).Value = 0.2501036404769519
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2113696666666667
$ws.Cells.Item(5, 14).Value = 0.634109
$ws.Cells.Item(5, 15).Value = 0.03795977003925348
$ws.Cells.Item(5, 16).Value = 0.03795977003925347
$ws.Cells.Item(5, 17).Value = 0.02826660643644445
$ws.Cells.Item(5, 18).Value = 0.254399457928
$ws.Cells.Item(5, 19).Value = 0.009493876678485225
$ws.Cells.Item(5, 20).Value = 0.009493876678485223

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col11a1"
$ws.Cells.Item(6, 3).Value = "Ddr1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.1337306666666667
$ws.Cells.Item(6, 8).Value = 0.401192
$ws.Cells.Item(6, 9).Value = 0.2501036404769519
$ws.Cells.Item(6, 10).Value = 0.2501036404769519
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.589504333333333
$ws.Cells.Item(6, 14).Value = 4.768513
$ws.Cells.Item(6, 15).Value = 0.2854582680725092
$ws.Cells.Item(6, 16).Value = 0.2854582680725091
$ws.Cells.Item(6, 17).Value = 0.2125654741662222
$ws.Cells.Item(6, 18).Value = 1.913089267496
$ws.Cells.Item(6, 19).Value = 0.0713941520491802
$ws.Cells.Item(6, 20).Value = 0.07139415204918019

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col11a1"
$ws.Cells.Item(7, 3).Value = "Ddr1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.1337306666666667
$ws.Cells.Item(7, 8).Value = 0.401192
$ws.Cells.Item(7, 9).Value = 0.2501036404769519
$ws.Cells.Item(7, 10).Value = 0.2501036404769519
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.767380666666666
$ws.Cells.Item(7, 14).Value = 11.302142
$ws.Cells.Item(7, 15).Value = 0.6765819618882374
$ws.Cells.Item(7, 16).Value = 0.6765819618882374
$ws.Cells.Item(7, 17).Value = 0.5038143281404445
$ws.Cells.Item(7, 18).Value = 4.534328953264
$ws.Cells.Item(7, 19).Value = 0.1692156117492865
$ws.Cells.Item(7, 20).Value = 0.1692156117492865

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Col11a1"
$ws.Cells.Item(8, 3).Value = "Ddr1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.3450066666666666
$ws.Cells.Item(8, 8).Value = 1.03502
$ws.Cells.Item(8, 9).Value = 0.6452328809309626
$ws.Cells.Item(8, 10).Value = 0.6452328809309626
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.2113696666666667
$ws.Cells.Item(8, 14).Value = 0.634109
$ws.Cells.Item(8, 15).Value = 0.03795977003925348
$ws.Cells.Item(8, 16).Value = 0.03795977003925347
$ws.Cells.Item(8, 17).Value = 0.07292394413111111
$ws.Cells.Item(8, 18).Value = 0.6563154971799999
$ws.Cells.Item(8, 19).Value = 0.02449289178190436
$ws.Cells.Item(8, 20).Value = 0.02449289178190436

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Col11a1"
$ws.Cells.Item(9, 3).Value = "Ddr1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.3450066666666666
$ws.Cells.Item(9, 8).Value = 1.03502
$ws.Cells.Item(9, 9).Value = 0.6452328809309626
$ws.Cells.Item(9, 10).Value = 0.6452328809309626
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.589504333333333
$ws.Cells.Item(9, 14).Value = 4.768513
$ws.Cells.Item(9, 15).Value = 0.2854582680725092
$ws.Cells.Item(9, 16).Value = 0.2854582680725091
$ws.Cells.Item(9, 17).Value = 0.5483895916955556
$ws.Cells.Item(9, 18).Value = 4.93550632526
$ws.Cells.Item(9, 19).Value = 0.1841870606939881
$ws.Cells.Item(9, 20).Value = 0.1841870606939881

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Col11a1"
$ws.Cells.Item(10, 3).Value = "Ddr1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.3450066666666666
$ws.Cells.Item(10, 8).Value = 1.03502
$ws.Cells.Item(10, 9).Value = 0.6452328809309626
$ws.Cells.Item(10, 10).Value = 0.6452328809309626
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.767380666666666
$ws.Cells.Item(10, 14).Value = 11.302142
$ws.Cells.Item(10, 15).Value = 0.6765819618882374
$ws.Cells.Item(10, 16).Value = 0.6765819618882374
$ws.Cells.Item(10, 17).Value = 1.299771445871111
$ws.Cells.Item(10, 18).Value = 11.69794301284
$ws.Cells.Item(10, 19).Value = 0.4365529284550702
$ws.Cells.Item(10, 20).Value = 0.4365529284550702

